$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "40.187.28") that must remain
# plain text exactly as in the source data (it is stored as inlineStr, not a number).
# Force text formatting before assigning so Excel does not auto-convert to a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.187.28"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.226.71"
$ws.Range("E3").Value = "  -0.69%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "295.28"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("E6").Value = "  +1.35%  "

$ws.Range("E7").Value = "  -1.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.479"
$ws.Range("E9").Value = "  +0.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.32"
$ws.Range("E10").Value = "  -3.42%  "

$ws.Range("E11").Value = "  -2.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "50.10"
$ws.Range("E12").Value = "  +6.09%  "

$ws.Range("E13").Value = "  +3.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.50"
$ws.Range("E14").Value = "  +0.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.578.21"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.84"
$ws.Range("E16").Value = "  -2.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.223.60"
$ws.Range("E17").Value = "  -0.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.734"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "40.105.15"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("E20").Value = "  -0.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.37"
$ws.Range("E21").Value = "  +4.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.82"
$ws.Range("E22").Value = "  -0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.62"
$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.43"
$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("E26").Value = "  -0.62%  "

$ws.Range("E27").Value = "  -2.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.82"
$ws.Range("E28").Value = "  -0.88%  "

$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.29"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.95"
$ws.Range("E31").Value = "  +1.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.15"
$ws.Range("E32").Value = "  -4.38%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.94"
$ws.Range("E34").Value = "  +0.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0718"
$ws.Range("E35").Value = "  -0.32%  "

$ws.Range("E36").Value = "  -1.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.91"
$ws.Range("E37").Value = "  +7.01%  "

$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.74"
$ws.Range("E39").Value = "  -5.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0977"
$ws.Range("E40").Value = "  -3.37%  "

$ws.Range("E41").Value = "  -1.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.138.96"
$ws.Range("E42").Value = "  +4.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.84"
$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0268"
$ws.Range("E46").Value = "  -1.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.75"
$ws.Range("E47").Value = "  -2.97%  "

$ws.Range("E48").Value = "  +4.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.442.48"
$ws.Range("E49").Value = "  -0.84%  "

$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.01"
$ws.Range("E51").Value = "  -0.90%  "

# Row 44 and 45: coin order swapped with updated data (EnergySwap now ranked above ApeXProtocol)
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.06"
$ws.Range("E44").Value = "  +10.38%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.12"
$ws.Range("E45").Value = "  -4.19%  "
